$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.402.37"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.458.45"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.67"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.49"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "2.454.89"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "2.900.78"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "62.168.94"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.452.65"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.53"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.96"
$ws.Range("E23").Value = "  -6.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.72"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.22"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "595.23"
$ws.Range("E27").Value = "  -5.39%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "0.0₃0961"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.92"
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.97"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.38"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.44"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.66"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.65"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.606"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  +13.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.83"
$ws.Range("E51").Value = "  -4.47%  "
